$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "66.649.02"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "3.233.43"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.232.86"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "3.756.64"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "66.669.35"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "3.227.99"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +28.35%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "506.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  +17.09%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.130"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.42%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "2.900.06"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.19%  "
